# relative entity was changed
# F2 previously contained an inappropriate remark about a student; replace it.
# A4 previously contained "S-E"; replace it with "NaN" (matching the other
# placeholder cells in column A).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("F2").Value = "kl,"
$ws.Range("A4").Value = "NaN"
